$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Invert the "divide by soma" ratio formulas in columns D and F ---
# Row 2 (plain, non-shared formulas)
$ws.Range("D2").Formula = "=C2/B2"
$ws.Range("F2").Formula = "=E2/B2"

# H2's formula (B2/G2) is removed entirely - cell becomes a blank, styled cell
$ws.Range("H2").ClearContents()

# Rows 5, 9, 10 participate in the shared formula group for D3:D11 / F3:F11.
# Re-point each populated member of the group to divide the other way round.
$ws.Range("D5").Formula  = "=C5/B5"
$ws.Range("D9").Formula  = "=C9/B9"
$ws.Range("D10").Formula = "=C10/B10"

$ws.Range("F5").Formula  = "=E5/B5"
$ws.Range("F9").Formula  = "=E9/B9"

# Row 11 loses its formula entirely (D11 and F11 become blank, styled cells)
$ws.Range("D11").ClearContents()
$ws.Range("F11").ClearContents()

# --- Update the saved selection / active cell ---
$ws.Range("G15").Select() | Out-Null
